$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 813
$wsExpo.Range("F4").Value = 279
$wsExpo.Range("F5").Value = 964
$wsExpo.Range("F6").Value = 2312
$wsExpo.Range("F7").Value = 198

# Sheet "全部类型" (all types): same events repeated, update column F too
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 813
$wsAll.Range("F4").Value = 279
$wsAll.Range("F7").Value = 964
$wsAll.Range("F8").Value = 2312
$wsAll.Range("F10").Value = 198
